$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Legs and Sesi 2 Update!" -- revised calibration readings for L1, L2, L3, R1, R2
$ws.Range("D4").Value = 1580
$ws.Range("F4").Value = 2010
$ws.Range("G4").Value = 2150

$ws.Range("C5").Value = 1570
$ws.Range("D5").Value = 1550

$ws.Range("C6").Value = 1350
$ws.Range("G6").Value = 2150

$ws.Range("C7").Value = 1410
$ws.Range("D7").Value = 1380

$ws.Range("C8").Value = 1600
$ws.Range("D8").Value = 1250

# Move the selection/active cell as it was left in the saved workbook
$ws.Activate()
$ws.Range("F17").Select()
